$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two new parameter rows (56, 57) use a style that is visually identical to the
# existing "body" rows that sit on the medium-blue fill with no border and left
# alignment (e.g. C23:D23) - font/fill/border/alignment all match, so we copy that
# formatting onto the new range before writing values.
$ws.Range("C23:D23").Copy()
$ws.Range("B56:D57").PasteSpecial(-4122)

# Row 56: MFTC_WEP_scaling parameter
$ws.Range("B56").Value = "MFTC_WEP_scaling"
# Write the numeric-looking value as literal text (quote-prefixed) so it is stored
# as a string, matching the source data - then re-apply the formatting on top so
# the quote-prefix doesn't leave behind a stray, differently-styled cell.
$ws.Range("C56").Value = "'1"
$ws.Range("D56").Value = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

# Row 57: WFF_or_Benefit parameter
$ws.Range("B57").Value = "WFF_or_Benefit"
$ws.Range("C57").Value = "Max"
$ws.Range("D57").Value = 'What work decision should we assume? Go off-benefit and receive IWTC = "WFF", stay on-benefit = "Benefit", or whichever gives a higher net income = "Max"'

# Re-apply formatting once more to keep the cell format consistent (quote-prefix
# entry above can otherwise leave C56 on a distinct auto-generated style).
$ws.Range("C23:D23").Copy()
$ws.Range("B56:D57").PasteSpecial(-4122)
